# Update "想去人数" (F column) values for several events that changed
# between the previous and the newly generated data snapshot.
#
# Sheet "展览" (1st sheet / sheetId 1):
#   F4:  213  -> 214
#   F6:  9805 -> 9808
#   F7:  884  -> 886
#   F10: 3390 -> 3892
#   F16: 535  -> 537
#   F19: 1432 -> 1437
#
# Sheet "全部类型" (4th sheet / sheetId 4) mirrors the same events one row lower:
#   F5:  213  -> 214
#   F7:  9805 -> 9808
#   F8:  884  -> 886
#   F11: 3390 -> 3892
#   F17: 535  -> 537
#   F20: 1432 -> 1437

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 214
$wsExhibition.Range("F6").Value  = 9808
$wsExhibition.Range("F7").Value  = 886
$wsExhibition.Range("F10").Value = 3892
$wsExhibition.Range("F16").Value = 537
$wsExhibition.Range("F19").Value = 1437

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 214
$wsAll.Range("F7").Value  = 9808
$wsAll.Range("F8").Value  = 886
$wsAll.Range("F11").Value = 3892
$wsAll.Range("F17").Value = 537
$wsAll.Range("F20").Value = 1437
